$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.048.00'
$ws.Cells.Item(2, 5).Value = '  +0.07%  '
$ws.Cells.Item(3, 4).Value = '1.910.92'
$ws.Cells.Item(3, 5).Value = '  +0.48%  '
$ws.Cells.Item(4, 4).Value = '''1.001'
$ws.Cells.Item(4, 5).Value = '  +0.06%  '
$ws.Cells.Item(5, 4).Value = '''0.7911'
$ws.Cells.Item(5, 5).Value = '  +6.62%  '
$ws.Cells.Item(6, 4).Value = '''243.01'
$ws.Cells.Item(6, 5).Value = '  -0.31%  '
$ws.Cells.Item(7, 5).Value = '  -0.05%  '
$ws.Cells.Item(8, 4).Value = '''0.3169'
$ws.Cells.Item(8, 5).Value = '  +3.28%  '
$ws.Cells.Item(9, 4).Value = '''26.38'
$ws.Cells.Item(9, 5).Value = '  +1.89%  '
$ws.Cells.Item(10, 4).Value = '''0.06929'
$ws.Cells.Item(10, 5).Value = '  +0.31%  '
$ws.Cells.Item(11, 4).Value = '''0.07986'
$ws.Cells.Item(12, 4).Value = '''0.7484'
$ws.Cells.Item(12, 5).Value = '  -1.68%  '
$ws.Cells.Item(13, 4).Value = '1.907.47'
$ws.Cells.Item(13, 5).Value = '  +0.26%  '
$ws.Cells.Item(14, 4).Value = '''5.237'
$ws.Cells.Item(14, 5).Value = '  -0.06%  '
$ws.Cells.Item(15, 4).Value = '''93.47'
$ws.Cells.Item(15, 5).Value = '  +2.31%  '
$ws.Cells.Item(16, 4).Value = '30.068.96'
$ws.Cells.Item(16, 5).Value = '  +0.11%  '
$ws.Cells.Item(17, 4).Value = '''14.05'
$ws.Cells.Item(17, 5).Value = '  +0.02%  '
$ws.Cells.Item(18, 4).Value = '''5.948'
$ws.Cells.Item(18, 5).Value = '  -4.46%  '
$ws.Cells.Item(19, 4).Value = '''246.98'
$ws.Cells.Item(19, 5).Value = '  +3.90%  '
$ws.Cells.Item(20, 4).Value = '''0.000007794'
$ws.Cells.Item(20, 5).Value = '  +0.49%  '
$ws.Cells.Item(21, 5).Value = '  +0.03%  '
$ws.Cells.Item(22, 5).Value = '  +0.04%  '
$ws.Cells.Item(23, 5).Value = '  -2.16%  '
$ws.Cells.Item(24, 4).Value = '''169.99'
$ws.Cells.Item(24, 5).Value = '  +2.14%  '
$ws.Cells.Item(25, 4).Value = '''9.317'
$ws.Cells.Item(25, 5).Value = '  +0.01%  '
$ws.Cells.Item(26, 4).Value = '''0.1391'
$ws.Cells.Item(26, 5).Value = '  +10.75%  '
$ws.Cells.Item(27, 4).Value = '''18.93'
$ws.Cells.Item(27, 5).Value = '  +0.59%  '
$ws.Cells.Item(28, 4).Value = '''2.046'
$ws.Cells.Item(28, 5).Value = '  +0.20%  '
$ws.Cells.Item(29, 4).Value = '''1.376'
$ws.Cells.Item(29, 5).Value = '  +1.89%  '
$ws.Cells.Item(30, 4).Value = '''1.523'
$ws.Cells.Item(30, 5).Value = '  -0.97%  '
$ws.Cells.Item(31, 4).Value = '''4.345'
$ws.Cells.Item(31, 5).Value = '  +1.01%  '
$ws.Cells.Item(32, 4).Value = '''0.05595'
$ws.Cells.Item(32, 5).Value = '  +5.41%  '
$ws.Cells.Item(33, 4).Value = '''4.116'
$ws.Cells.Item(33, 5).Value = '  +1.63%  '
$ws.Cells.Item(34, 5).Value = '  -2.24%  '
$ws.Cells.Item(35, 4).Value = '''0.7376'
$ws.Cells.Item(35, 5).Value = '  -0.28%  '
$ws.Cells.Item(36, 4).Value = '''2.730'
$ws.Cells.Item(36, 5).Value = '  +0.15%  '
$ws.Cells.Item(37, 4).Value = '''0.01941'
$ws.Cells.Item(37, 5).Value = '  +0.01%  '
$ws.Cells.Item(38, 4).Value = '''2.794'
$ws.Cells.Item(38, 5).Value = '  +0.57%  '
$ws.Cells.Item(39, 4).Value = '''6.189'
$ws.Cells.Item(39, 5).Value = '  -1.40%  '
$ws.Cells.Item(40, 4).Value = '''0.4445'
$ws.Cells.Item(40, 5).Value = '  -0.12%  '
$ws.Cells.Item(41, 4).Value = '''72.67'
$ws.Cells.Item(41, 5).Value = '  -0.66%  '
$ws.Cells.Item(42, 4).Value = '''1.001'
$ws.Cells.Item(42, 5).Value = '  +0.02%  '
$ws.Cells.Item(43, 5).Value = '  -3.22%  '
$ws.Cells.Item(44, 4).Value = '''0.8351'
$ws.Cells.Item(44, 5).Value = '  -0.23%  '
$ws.Cells.Item(45, 4).Value = '''7.573'
$ws.Cells.Item(45, 5).Value = '  -0.79%  '
$ws.Cells.Item(46, 4).Value = '''100.59'
$ws.Cells.Item(46, 5).Value = '  -0.70%  '
$ws.Cells.Item(47, 4).Value = '''9.811'
$ws.Cells.Item(47, 5).Value = '  -0.04%  '
$ws.Cells.Item(48, 4).Value = '''989.08'
$ws.Cells.Item(48, 5).Value = '  +8.33%  '
$ws.Cells.Item(49, 4).Value = '2.064.61'
$ws.Cells.Item(49, 5).Value = '  +0.66%  '
$ws.Cells.Item(50, 4).Value = '''36.32'
$ws.Cells.Item(50, 5).Value = '  -1.14%  '
$ws.Cells.Item(51, 5).Value = '  +2.95%  '
